# Update "want to go" counts (column F) on several sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12016
$ws1.Range("F4").Value = 35
$ws1.Range("F8").Value = 11911
$ws1.Range("F10").Value = 1177
$ws1.Range("F11").Value = 108
$ws1.Range("F13").Value = 1792
$ws1.Range("F14").Value = 5897
$ws1.Range("F16").Value = 3551
$ws1.Range("F17").Value = 195

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 576

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 576
$ws4.Range("F5").Value = 12016
$ws4.Range("F6").Value = 35
$ws4.Range("F11").Value = 11911
$ws4.Range("F13").Value = 1177
$ws4.Range("F14").Value = 108
$ws4.Range("F16").Value = 1792
$ws4.Range("F18").Value = 5897
$ws4.Range("F20").Value = 3551
$ws4.Range("F21").Value = 195
